$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = @'
MATCH (ss:study_subject)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sample_of_study_subject]-(samp:sample)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)

WITH DISTINCT ss, samp, collect(DISTINCT samp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"]  and d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"] and f.file_format IN ["txt"]
return DISTINCT ss.study_subject_id as `Case ID`,
   p.program_acronym as `Program Code`,
    p.program_id as Program_ID,
   s.study_acronym as `Arm`,
   ss.disease_subtype as `Diagnosis`,
   sf.grouped_recurrence_score AS `Recurrence Score`,
   d.tumor_size_group AS `tumor_size`,
   d.er_status AS `ER Status`,
   d.pr_status AS `PR Status`,
   demo.age_at_index AS `Age (years)`,
	demo.survival_time AS `Survival (days)`
'@
$c2 = @'
MATCH (ss:study_subject)
MATCH (samp:sample)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"]  and d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"] and f.file_format IN ["txt"]
WITH ss 
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@
$b3 = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"] and d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"]and f.file_format IN ["txt"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`,
            lp.test_name as Platform
'@
$c3 = @'
MATCH (ss:study_subject)
MATCH (samp:sample)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"]  and d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"] and f.file_format IN ["txt"]
WITH ss 
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@
$b4 = @'
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"] and  d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"]  and f.file_format IN ["txt"]
RETURN  f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
    f.file_size AS `Size`,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`,
    ss.disease_subtype as `Diagnosis`
'@
$c4 = @'
MATCH (ss:study_subject)
MATCH (samp:sample)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
WHERE ss.disease_subtype IN ["Tubular Carcinoma"] and d.tumor_size_group In ["(3,4]"]  and d.pr_status In ["Positive"] and tp.endocrine_therapy_type IN ["Tam & AI"] and f.file_format IN ["txt"]
WITH ss 
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@

$ws.Range("C2").Value2 = $c2
$ws.Range("B2").Value2 = $b2
$ws.Range("B3").Value2 = $b3
$ws.Range("C3").Value2 = $c3
$ws.Range("B4").Value2 = $b4
$ws.Range("C4").Value2 = $c4

$ws.Range("C4").Select()
